# Daily scrape update - 2026-01-01 03:45:33 UTC
# Updates the opportunities sheet: refreshes rows 2-6 with new listings
# (all from Wavetec / Karachi, Pakistan) and appends a new row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (characters) ---
# Input ColumnWidth values are offset by -0.875 from the desired stored
# width because of Excel's internal character-width padding/rounding.
$ws.Columns.Item(3).ColumnWidth = 30.125   # C: 53 -> 31
$ws.Columns.Item(4).ColumnWidth = 19.125   # D: 36 -> 20
$ws.Columns.Item(6).ColumnWidth = 16.125   # F: 16 -> 17
$ws.Columns.Item(8).ColumnWidth = 24.125   # H: 21 -> 25

# --- Keep the OPPORTUNITY ID column formatted as text so the numeric-looking
#     IDs are stored as strings rather than being coerced into numbers ---
$ws.Range("A2:A7").NumberFormat = "@"

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = "1323507"
$ws.Cells.Item(2, 2).Value = "https://aiesec.org/opportunity/global-talent/1323507"
$ws.Cells.Item(2, 3).Value = "Sales Intern"
$ws.Cells.Item(2, 4).Value = "Karachi, Pakistan"
$ws.Cells.Item(2, 5).Value = "No"
$ws.Cells.Item(2, 6).Value = "14 applicants"
$ws.Cells.Item(2, 7).Value = "6 - 18 Months"
$ws.Cells.Item(2, 8).Value = "Wavetec"

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = "1323506"
$ws.Cells.Item(3, 2).Value = "https://aiesec.org/opportunity/global-talent/1323506"
$ws.Cells.Item(3, 3).Value = "Marketing Intern"
$ws.Cells.Item(3, 4).Value = "Karachi, Pakistan"
$ws.Cells.Item(3, 5).Value = "No"
$ws.Cells.Item(3, 6).Value = "5 applicants"
$ws.Cells.Item(3, 7).Value = "6 - 18 Months"
$ws.Cells.Item(3, 8).Value = "Wavetec"

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value = "1323504"
$ws.Cells.Item(4, 2).Value = "https://aiesec.org/opportunity/global-talent/1323504"
$ws.Cells.Item(4, 3).Value = "Management Trainee Officer"
$ws.Cells.Item(4, 4).Value = "Karachi, Pakistan"
$ws.Cells.Item(4, 5).Value = "No"
$ws.Cells.Item(4, 6).Value = "16 applicants"
$ws.Cells.Item(4, 7).Value = "6 - 18 Months"
$ws.Cells.Item(4, 8).Value = "Wavetec"

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value = "1323503"
$ws.Cells.Item(5, 2).Value = "https://aiesec.org/opportunity/global-talent/1323503"
$ws.Cells.Item(5, 3).Value = "Customer Support Engineer"
$ws.Cells.Item(5, 4).Value = "Karachi, Pakistan"
$ws.Cells.Item(5, 5).Value = "No"
$ws.Cells.Item(5, 6).Value = "7 applicants"
$ws.Cells.Item(5, 7).Value = "6 - 18 Months"
$ws.Cells.Item(5, 8).Value = "Wavetec"

# --- Row 6 ---
$ws.Cells.Item(6, 1).Value = "1323391"
$ws.Cells.Item(6, 2).Value = "https://aiesec.org/opportunity/global-talent/1323391"
$ws.Cells.Item(6, 3).Value = "Business Development Officer"
$ws.Cells.Item(6, 4).Value = "Karachi, Pakistan"
$ws.Cells.Item(6, 5).Value = "No"
$ws.Cells.Item(6, 6).Value = "9 applicants"
$ws.Cells.Item(6, 7).Value = "6 - 18 Months"
$ws.Cells.Item(6, 8).Value = "Wavetec"

# --- Row 7 (new) ---
$ws.Cells.Item(7, 1).Value = "1309138"
$ws.Cells.Item(7, 2).Value = "https://aiesec.org/opportunity/global-talent/1309138"
$ws.Cells.Item(7, 3).Value = "Business Development"
$ws.Cells.Item(7, 4).Value = "İstanbul, Türkiye"
$ws.Cells.Item(7, 5).Value = "No"
$ws.Cells.Item(7, 6).Value = "113 applicants"
$ws.Cells.Item(7, 7).Value = "6 - 18 Months"
$ws.Cells.Item(7, 8).Value = "Dentekay Dental Clinic"
